$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Affected Flight Number" column (column A) is no longer part of the
# manifest upload template, so remove it entirely. Excel shifts every
# column after it one position to the left, which also removes the now
# unused "Affected Flight Number" shared string and renumbers the
# worksheet dimension/columns accordingly.
$ws.Columns.Item(1).Delete() | Out-Null

# After the delete, reselect near the top of the sheet (mirrors the
# author's subsequent click into B2 while continuing the manifest work).
$ws.Range("B2").Select() | Out-Null

$wb.Save()
